$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 87.875
$ws.Range("I9").Value = 106.833336
$ws.Range("K9").Value = 106.833336
$ws.Range("M9").Value = 62.166664
$ws.Range("H11").Value = 161.94118
$ws.Range("I11").Value = 161.94118
$ws.Range("K11").Value = 161.94118
$ws.Range("M11").Value = -21.94118
$ws.Range("H18").Value = 2136.2727
$ws.Range("I18").Value = 2136.2727
$ws.Range("K18").Value = 2136.2727
$ws.Range("M18").Value = -1852.2727
$ws.Range("H76").Value = 2675.5
$ws.Range("I76").Value = 1966.6666
$ws.Range("K76").Value = 1966.6666
$ws.Range("M76").Value = -1651.6666
$ws.Range("H79").Value = 2675.5
$ws.Range("I79").Value = 1966.6666
$ws.Range("K79").Value = 1966.6666
$ws.Range("M79").Value = -874.6666
$ws.Range("H132").Value = 1089.1936
$ws.Range("I132").Value = 792.23334
$ws.Range("K132").Value = 2376.70002
$ws.Range("M132").Value = 153.2999799999998

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 0
$ws.Range("I2").Value = 0
$ws.Range("K2").Value = 0
$ws.Range("M2").Value = ""
$ws.Range("H61").Value = 1414.4
$ws.Range("I61").Value = 1414.4
$ws.Range("K61").Value = 1414.4
$ws.Range("M61").Value = -1202.4
$ws.Range("H74").Value = 3052
$ws.Range("I74").Value = 1153
$ws.Range("J74").Value = 6850
$ws.Range("K74").Value = 1153
$ws.Range("L74").Value = 6850
$ws.Range("M74").Value = -279
$ws.Range("N74").Value = -8598
$ws.Range("H77").Value = 3052
$ws.Range("I77").Value = 1153
$ws.Range("J77").Value = 6850
$ws.Range("K77").Value = 5765
$ws.Range("L77").Value = 34250
$ws.Range("M77").Value = -1397
$ws.Range("N77").Value = -42986
$ws.Range("H80").Value = 89998.5
$ws.Range("J80").Value = 89998.5
$ws.Range("L80").Value = 89998.5
$ws.Range("N80").Value = -91994.5
$ws.Range("H83").Value = 89998.5
$ws.Range("J83").Value = 89998.5
$ws.Range("L83").Value = 269995.5
$ws.Range("N83").Value = -279979.5
$ws.Range("H97").Value = 2247.8
$ws.Range("I97").Value = 2690
$ws.Range("K97").Value = 2690
$ws.Range("M97").Value = -2194
$ws.Range("H116").Value = 0
$ws.Range("I116").Value = 0
$ws.Range("K116").Value = 0
$ws.Range("M116").Value = ""
$ws.Range("H119").Value = 0
$ws.Range("J119").Value = 0
$ws.Range("L119").Value = 0
$ws.Range("N119").Value = ""
$ws.Range("H132").Value = 0
$ws.Range("I132").Value = 0
$ws.Range("K132").Value = 0
$ws.Range("M132").Value = ""
$ws.Range("H136").Value = 1414.4
$ws.Range("I136").Value = 1414.4
$ws.Range("K136").Value = 4243.200000000001
$ws.Range("M136").Value = -1693.200000000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 0
$ws.Range("I3").Value = 0
$ws.Range("K3").Value = 0
$ws.Range("M3").Value = ""
$ws.Range("H80").Value = 491
$ws.Range("I80").Value = 1450
$ws.Range("J80").Value = 299.2
$ws.Range("K80").Value = 1450
$ws.Range("L80").Value = 299.2
$ws.Range("M80").Value = -452
$ws.Range("N80").Value = -2295.2
$ws.Range("H83").Value = 491
$ws.Range("I83").Value = 1450
$ws.Range("J83").Value = 299.2
$ws.Range("K83").Value = 7250
$ws.Range("L83").Value = 1496
$ws.Range("M83").Value = -2258
$ws.Range("N83").Value = -11480
$ws.Range("H86").Value = 2224.25
$ws.Range("I86").Value = 2448
$ws.Range("K86").Value = 2448
$ws.Range("M86").Value = -1325
$ws.Range("H89").Value = 2224.25
$ws.Range("I89").Value = 2448
$ws.Range("K89").Value = 12240
$ws.Range("M89").Value = -6624
$ws.Range("H99").Value = 2644.2778
$ws.Range("I99").Value = 2666.6667
$ws.Range("K99").Value = 2666.6667
$ws.Range("M99").Value = -1168.6667
$ws.Range("H128").Value = 0
$ws.Range("I128").Value = 0
$ws.Range("K128").Value = 0
$ws.Range("M128").Value = ""
$ws.Range("H134").Value = 2587.111
$ws.Range("I134").Value = 2587.111
$ws.Range("K134").Value = 7761.333
$ws.Range("M134").Value = -5226.333

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 2366.9167
$ws.Range("I58").Value = 1266.1111
$ws.Range("K58").Value = 1266.1111
$ws.Range("M58").Value = -1063.1111
$ws.Range("H97").Value = 0
$ws.Range("J97").Value = 0
$ws.Range("L97").Value = 0
$ws.Range("N97").Value = ""
$ws.Range("H99").Value = 16703.354
$ws.Range("I99").Value = 14163.5
$ws.Range("J99").Value = 18088.727
$ws.Range("K99").Value = 14163.5
$ws.Range("L99").Value = 18088.727
$ws.Range("M99").Value = -12665.5
$ws.Range("N99").Value = -21084.727
$ws.Range("H109").Value = 56129.5
$ws.Range("J109").Value = 62000
$ws.Range("L109").Value = 62000
$ws.Range("N109").Value = -64080
$ws.Range("H126").Value = 16703.354
$ws.Range("I126").Value = 14163.5
$ws.Range("J126").Value = 18088.727
$ws.Range("K126").Value = 42490.5
$ws.Range("L126").Value = 54266.181
$ws.Range("M126").Value = -40020.5
$ws.Range("N126").Value = -59206.181
$ws.Range("H132").Value = 1464.5834
$ws.Range("I132").Value = 1523
$ws.Range("J132").Value = 1242.6
$ws.Range("K132").Value = 4569
$ws.Range("L132").Value = 3727.8
$ws.Range("M132").Value = -2039
$ws.Range("N132").Value = -8787.799999999999
$ws.Range("H134").Value = 2311.7097
$ws.Range("I134").Value = 2136.4783
$ws.Range("J134").Value = 2815.5
$ws.Range("K134").Value = 6409.4349
$ws.Range("L134").Value = 8446.5
$ws.Range("M134").Value = -3874.4349
$ws.Range("N134").Value = -13516.5
$ws.Range("H136").Value = 2366.9167
$ws.Range("I136").Value = 1266.1111
$ws.Range("K136").Value = 3798.3333
$ws.Range("M136").Value = -1248.3333

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H10").Value = 59.25
$ws.Range("I10").Value = 64.666664
$ws.Range("J10").Value = 43
$ws.Range("K10").Value = 193.999992
$ws.Range("L10").Value = 129
$ws.Range("M10").Value = -54.99999199999999
$ws.Range("N10").Value = -407
$ws.Range("H11").Value = 86154.42999999999
$ws.Range("I11").Value = 150020.25
$ws.Range("J11").Value = 1000
$ws.Range("K11").Value = 450060.75
$ws.Range("L11").Value = 3000
$ws.Range("M11").Value = -449920.75
$ws.Range("N11").Value = -3280
$ws.Range("H35").Value = 300
$ws.Range("I35").Value = 0
$ws.Range("J35").Value = 300
$ws.Range("K35").Value = 0
$ws.Range("L35").Value = 900
$ws.Range("M35").Value = ""
$ws.Range("N35").Value = -1476
$ws.Range("H46").Value = 2500592
$ws.Range("J46").Value = 3333872.8
$ws.Range("L46").Value = 10001618.4
$ws.Range("N46").Value = -10001800.4
$ws.Range("H51").Value = 2575
$ws.Range("I51").Value = 2362.5
$ws.Range("K51").Value = 7087.5
$ws.Range("M51").Value = -6627.5
$ws.Range("H107").Value = 2198.0833
$ws.Range("I107").Value = 395.8
$ws.Range("J107").Value = 3485.4285
$ws.Range("K107").Value = 1187.4
$ws.Range("L107").Value = 10456.2855
$ws.Range("M107").Value = 732.5999999999999
$ws.Range("N107").Value = -14296.2855
$ws.Range("H114").Value = 0
$ws.Range("I114").Value = 0
$ws.Range("J114").Value = 0
$ws.Range("K114").Value = 0
$ws.Range("L114").Value = 0
$ws.Range("M114").Value = ""
$ws.Range("N114").Value = ""
$ws.Range("H118").Value = 1162.2
$ws.Range("I118").Value = 897.25
$ws.Range("K118").Value = 2691.75
$ws.Range("M118").Value = -1448.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 3328.2856
$ws.Range("I126").Value = 2321.125
$ws.Range("K126").Value = 6963.375
$ws.Range("M126").Value = -4493.375

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5275
$ws.Range("I7").Value = 5275
$ws.Range("K7").Value = 5275
$ws.Range("M7").Value = -5163
$ws.Range("H40").Value = 3099.5386
$ws.Range("I40").Value = 2982.9167
$ws.Range("K40").Value = 2982.9167
$ws.Range("M40").Value = -2846.9167
$ws.Range("H116").Value = 99680
$ws.Range("J116").Value = 99680
$ws.Range("L116").Value = 99680
$ws.Range("N116").Value = -108858
$ws.Range("H123").Value = 35000
$ws.Range("I123").Value = 35000
$ws.Range("K123").Value = 35000
$ws.Range("M123").Value = -30100
$ws.Range("H126").Value = 5275
$ws.Range("I126").Value = 5275
$ws.Range("K126").Value = 15825
$ws.Range("M126").Value = -13355

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H7").Value = 2249.75
$ws.Range("I7").Value = 2999
$ws.Range("J7").Value = 2000
$ws.Range("K7").Value = 2999
$ws.Range("L7").Value = 2000
$ws.Range("M7").Value = -2886
$ws.Range("N7").Value = -2226
$ws.Range("H81").Value = 15571.286
$ws.Range("I81").Value = 19999
$ws.Range("J81").Value = 14833.333
$ws.Range("K81").Value = 39998
$ws.Range("L81").Value = 29666.666
$ws.Range("M81").Value = -38937
$ws.Range("N81").Value = -31788.666
$ws.Range("H84").Value = 15571.286
$ws.Range("I84").Value = 19999
$ws.Range("J84").Value = 14833.333
$ws.Range("K84").Value = 199990
$ws.Range("L84").Value = 148333.33
$ws.Range("M84").Value = -194686
$ws.Range("N84").Value = -158941.33
$ws.Range("H131").Value = 0
$ws.Range("J131").Value = 0
$ws.Range("L131").Value = 0
$ws.Range("N131").Value = ""

Write-Output "applied market data refresh"